$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 9: #property=8, address=高松市 (same as row8), date=2024.3.31, numberOfInsured=68290
$ws.Range("A9").Value = 8
$ws.Range("A9").HorizontalAlignment = -4131

$ws.Range("B9").Value = $ws.Range("B8").Value2

$ws.Range("C9").Formula = "=""2024.3.31"""
$ws.Range("C9").Copy()
$ws.Range("C9").PasteSpecial(-4163)

$ws.Range("D9").Value = 68290
$ws.Range("D9").HorizontalAlignment = -4131

$null = $ws.Range("G10").Select()
